# Applies "new sensitivity and calculus" changes:
#  - Sheet 1 "Model Accuracy": add Market threshold/min/max, Recall, Precision
#    columns (C:G) and update Accuracy (col B) values.
#  - Sheets 2-6 "Confusion Matrix *": update the confusion-matrix counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Model Accuracy (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cells C1:G1 - copy the existing bold/bordered header style from
# B1 and then overwrite the values with the new header text.
$ws1.Range("B1").Copy($ws1.Range("C1:G1"))
$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("D1").Value = "Market min"
$ws1.Range("E1").Value = "Market max"
$ws1.Range("F1").Value = "Recall"
$ws1.Range("G1").Value = "Precision"

# Row 2 - TOTALENERGIES SE
$ws1.Range("B2").Value = 50.42787286063569
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 11.11111111111111
$ws1.Range("G2").Value = 1.538461538461539

# Row 3 - FMC CORP
$ws1.Range("B3").Value = 30.0122249388753
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 10.45576407506702
$ws1.Range("G3").Value = 24.07407407407407

# Row 4 - BP PLC
$ws1.Range("B4").Value = 75.61124694376528
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

# Row 5 - STORA ENSO
$ws1.Range("B5").Value = 61.7359413202934
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 5.454545454545454
$ws1.Range("G5").Value = 5.263157894736842

# Row 6 - BHP GROUP
$ws1.Range("B6").Value = 73.59413202933986
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 25
$ws1.Range("G6").Value = 0.411522633744856

# ---------------------------------------------------------------------------
# Sheet 2: Confusion Matrix TOTALENERGIES SE (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = 818
$ws2.Range("D3").Value = 4

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix FMC CORP (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 39
$ws3.Range("C2").Value = 89
$ws3.Range("D2").Value = 34

$ws3.Range("B3").Value = 181
$ws3.Range("C3").Value = 310
$ws3.Range("D3").Value = 180

$ws3.Range("B4").Value = 153
$ws3.Range("C4").Value = 253
$ws3.Range("D4").Value = 142

# ---------------------------------------------------------------------------
# Sheet 4: Confusion Matrix BP PLC (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B3").Value = 29
$ws4.Range("C3").Value = 1226
$ws4.Range("D3").Value = 31

$ws4.Range("B4").Value = 11
$ws4.Range("C4").Value = 283
$ws4.Range("D4").Value = 11

# ---------------------------------------------------------------------------
# Sheet 5: Confusion Matrix STORA ENSO (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = 6
$ws5.Range("C2").Value = 103
$ws5.Range("D2").Value = 5

$ws5.Range("B3").Value = 81
$ws5.Range("C3").Value = 985
$ws5.Range("D3").Value = 83

$ws5.Range("B4").Value = 23
$ws5.Range("C4").Value = 270
$ws5.Range("D4").Value = 19

# ---------------------------------------------------------------------------
# Sheet 6: Confusion Matrix BHP GROUP (-0.2, 0.2, 0.2)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B2").Value = 1
$ws6.Range("C2").Value = 241
$ws6.Range("D2").Value = 1

$ws6.Range("B3").Value = 3
$ws6.Range("C3").Value = 1203
$ws6.Range("D3").Value = 2

$ws6.Range("B4").Value = 0
$ws6.Range("C4").Value = 129
$ws6.Range("D4").Value = 0
